$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = 1
$ws.Range("B10").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C10").Value = "Arica y Parinacota"
$ws.Range("D10").Value = Get-Date -Year 2022 -Month 1 -Day 7 -Hour 0 -Minute 0 -Second 0
$ws.Range("D10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E10").Value = 15
$ws.Range("F10").Value = "Fruta"
$ws.Range("G10").Value = 100103
$ws.Range("H10").Value = "Frutos de hueso (carozo)"
$ws.Range("I10").Value = 100103001
$ws.Range("J10").Value = "Cereza"
$ws.Range("K10").Value = "Santina"
$ws.Range("L10").Value = "Segunda"
$ws.Range("M10").Value = 200
$ws.Range("N10").Value = 15000
$ws.Range("O10").Value = 16000
$ws.Range("P10").Value = 15500
$ws.Range("Q10").Value = "$/bandeja 12 kilos"
$ws.Range("R10").Value = "Región de O'Higgins"
$ws.Range("S10").Value = 1292
$ws.Range("T10").Value = 12
